$wb = $excel.ActiveWorkbook

# --- Sheet Caso1: add Autogluon / H2O / AutoSklearn columns ---
$ws = $wb.Worksheets.Item("Caso1")

# Header row: new column titles (style matches existing unstyled G1 header cell)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("H1").Style = $ws.Range("G1").Style
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("I1").Style = $ws.Range("G1").Style
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("J1").Style = $ws.Range("G1").Style

$ws.Range("G2").Value = 0.9762083292007446
$ws.Range("H2").Value = 0.9727476835250854
$ws.Range("H2").Style = $ws.Range("G2").Style
$ws.Range("I2").Value = 0.9732129894835398
$ws.Range("I2").Style = $ws.Range("G2").Style
$ws.Range("J2").Value = 0.972614660859108
$ws.Range("J2").Style = $ws.Range("G2").Style
$ws.Range("G3").Value = 0.9704018831253052
$ws.Range("H3").Value = 0.9706719517707825
$ws.Range("H3").Style = $ws.Range("G3").Style
$ws.Range("I3").Value = 0.9716591719662176
$ws.Range("I3").Style = $ws.Range("G3").Style
$ws.Range("J3").Value = 0.9711216222494841
$ws.Range("J3").Style = $ws.Range("G3").Style
$ws.Range("G4").Value = 0.9717544913291931
$ws.Range("H4").Value = 0.9700327515602112
$ws.Range("H4").Style = $ws.Range("G4").Style
$ws.Range("I4").Value = 0.9713593756214742
$ws.Range("I4").Style = $ws.Range("G4").Style
$ws.Range("J4").Value = 0.9708215743303299
$ws.Range("J4").Style = $ws.Range("G4").Style
$ws.Range("G5").Value = 0.9701979160308838
$ws.Range("H5").Value = 0.9698159694671631
$ws.Range("H5").Style = $ws.Range("G5").Style
$ws.Range("I5").Value = 0.9713398369162124
$ws.Range("I5").Style = $ws.Range("G5").Style
$ws.Range("J5").Value = 0.9710562583059072
$ws.Range("J5").Style = $ws.Range("G5").Style
$ws.Range("G6").Value = 0.9702368378639221
$ws.Range("H6").Value = 0.9697266817092896
$ws.Range("H6").Style = $ws.Range("G6").Style
$ws.Range("I6").Value = 0.9712923521118556
$ws.Range("I6").Style = $ws.Range("G6").Style
$ws.Range("J6").Value = 0.9707390088587999
$ws.Range("J6").Style = $ws.Range("G6").Style
$ws.Range("G7").Value = 0.9681689739227295
$ws.Range("H7").Value = 0.9699131846427917
$ws.Range("H7").Style = $ws.Range("G7").Style
$ws.Range("I7").Value = 0.971309960266814
$ws.Range("I7").Style = $ws.Range("G7").Style
$ws.Range("J7").Value = 0.9708166979253292
$ws.Range("J7").Style = $ws.Range("G7").Style
$ws.Range("G8").Value = 0.9688431620597839
$ws.Range("H8").Value = 0.9698947668075562
$ws.Range("H8").Style = $ws.Range("G8").Style
$ws.Range("I8").Value = 0.971337542150144
$ws.Range("I8").Style = $ws.Range("G8").Style
$ws.Range("J8").Value = 0.9707022868096828
$ws.Range("J8").Style = $ws.Range("G8").Style
$ws.Range("G9").Value = 0.9698131680488586
$ws.Range("H9").Value = 0.9697583913803101
$ws.Range("H9").Style = $ws.Range("G9").Style
$ws.Range("I9").Value = 0.9712673561827564
$ws.Range("I9").Style = $ws.Range("G9").Style
$ws.Range("J9").Value = 0.9708121549338102
$ws.Range("J9").Style = $ws.Range("G9").Style
$ws.Range("G10").Value = 0.9690741896629333
$ws.Range("H10").Value = 0.9699063301086426
$ws.Range("H10").Style = $ws.Range("G10").Style
$ws.Range("I10").Value = 0.9712177561538838
$ws.Range("I10").Style = $ws.Range("G10").Style
$ws.Range("J10").Value = 0.9706190526485443
$ws.Range("J10").Style = $ws.Range("G10").Style
$ws.Range("G11").Value = 0.9688156843185425
$ws.Range("H11").Value = 0.9701405167579651
$ws.Range("H11").Style = $ws.Range("G11").Style
$ws.Range("I11").Value = 0.9712176087838172
$ws.Range("I11").Style = $ws.Range("G11").Style
$ws.Range("J11").Value = 0.9705638475716114
$ws.Range("J11").Style = $ws.Range("G11").Style
$ws.Range("G12").Value = 0.9693359732627869
$ws.Range("H12").Value = 0.9699406027793884
$ws.Range("H12").Style = $ws.Range("G12").Style
$ws.Range("I12").Value = 0.9712173785535176
$ws.Range("I12").Style = $ws.Range("G12").Style
$ws.Range("J12").Value = 0.9706821534782648
$ws.Range("J12").Style = $ws.Range("G12").Style
$ws.Range("G13").Value = 0.9678239822387695
$ws.Range("H13").Value = 0.969895601272583
$ws.Range("H13").Style = $ws.Range("G13").Style
$ws.Range("I13").Value = 0.9712225517584312
$ws.Range("I13").Style = $ws.Range("G13").Style
$ws.Range("J13").Value = 0.9706613644957542
$ws.Range("J13").Style = $ws.Range("G13").Style
$ws.Range("G14").Value = 0.9709464311599731
$ws.Range("H14").Value = 0.9702876806259155
$ws.Range("H14").Style = $ws.Range("G14").Style
$ws.Range("I14").Value = 0.9714898582650204
$ws.Range("I14").Style = $ws.Range("G14").Style
$ws.Range("J14").Value = 0.9710502363741398
$ws.Range("J14").Style = $ws.Range("G14").Style
$ws.Range("G15").Value = 0.9726852774620056
$ws.Range("H15").Value = 0.9703362584114075
$ws.Range("H15").Style = $ws.Range("G15").Style
$ws.Range("I15").Value = 0.9715341641561924
$ws.Range("I15").Style = $ws.Range("G15").Style
$ws.Range("J15").Value = 0.9711104389280081
$ws.Range("J15").Style = $ws.Range("G15").Style
$ws.Range("G16").Value = 0.969898521900177
$ws.Range("H16").Value = 0.9699578285217285
$ws.Range("H16").Style = $ws.Range("G16").Style
$ws.Range("I16").Value = 0.9713296706054768
$ws.Range("I16").Style = $ws.Range("G16").Style
$ws.Range("J16").Value = 0.9707029983401299
$ws.Range("J16").Style = $ws.Range("G16").Style
$ws.Range("G17").Value = 0.9712648391723633
$ws.Range("H17").Value = 0.9701982140541077
$ws.Range("H17").Style = $ws.Range("G17").Style
$ws.Range("I17").Value = 0.9713044433872025
$ws.Range("I17").Style = $ws.Range("G17").Style
$ws.Range("J17").Value = 0.9706152696162462
$ws.Range("J17").Style = $ws.Range("G17").Style
$ws.Range("G18").Value = 0.9708113074302673
$ws.Range("H18").Value = 0.9698747992515564
$ws.Range("H18").Style = $ws.Range("G18").Style
$ws.Range("I18").Value = 0.9713279140427258
$ws.Range("I18").Style = $ws.Range("G18").Style
$ws.Range("J18").Value = 0.9706277046352625
$ws.Range("J18").Style = $ws.Range("G18").Style
$ws.Range("G19").Value = 0.9702373147010803
$ws.Range("H19").Value = 0.9699745178222656
$ws.Range("H19").Style = $ws.Range("G19").Style
$ws.Range("I19").Value = 0.971320653063784
$ws.Range("I19").Style = $ws.Range("G19").Style
$ws.Range("J19").Value = 0.9706914722919464
$ws.Range("J19").Style = $ws.Range("G19").Style

# --- Sheet Caso2: add Autogluon / H2O / AutoSklearn columns ---
$ws = $wb.Worksheets.Item("Caso2")

# Header row: new column titles (style matches existing unstyled G1 header cell)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("H1").Style = $ws.Range("G1").Style
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("I1").Style = $ws.Range("G1").Style
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("J1").Style = $ws.Range("G1").Style

$ws.Range("G2").Value = 0.9765862226486206
$ws.Range("H2").Value = 0.9729588627815247
$ws.Range("H2").Style = $ws.Range("G2").Style
$ws.Range("I2").Value = 0.9733498185746912
$ws.Range("I2").Style = $ws.Range("G2").Style
$ws.Range("J2").Value = 0.9728307984769344
$ws.Range("J2").Style = $ws.Range("G2").Style
$ws.Range("G3").Value = 0.9708735942840576
$ws.Range("H3").Value = 0.9709711074829102
$ws.Range("H3").Style = $ws.Range("G3").Style
$ws.Range("I3").Value = 0.9718379325155314
$ws.Range("I3").Style = $ws.Range("G3").Style
$ws.Range("J3").Value = 0.9714913461357355
$ws.Range("J3").Style = $ws.Range("G3").Style
$ws.Range("G4").Value = 0.9722869396209717
$ws.Range("H4").Value = 0.97056645154953
$ws.Range("H4").Style = $ws.Range("G4").Style
$ws.Range("I4").Value = 0.971539040624009
$ws.Range("I4").Style = $ws.Range("G4").Style
$ws.Range("J4").Value = 0.9711607340723276
$ws.Range("J4").Style = $ws.Range("G4").Style
$ws.Range("G5").Value = 0.9706631898880005
$ws.Range("H5").Value = 0.9703032374382019
$ws.Range("H5").Style = $ws.Range("G5").Style
$ws.Range("I5").Value = 0.9715200697457376
$ws.Range("I5").Style = $ws.Range("G5").Style
$ws.Range("J5").Value = 0.9712864831089973
$ws.Range("J5").Style = $ws.Range("G5").Style
$ws.Range("G6").Value = 0.9707242846488953
$ws.Range("H6").Value = 0.9701356291770935
$ws.Range("H6").Style = $ws.Range("G6").Style
$ws.Range("I6").Value = 0.971473245815874
$ws.Range("I6").Style = $ws.Range("G6").Style
$ws.Range("J6").Value = 0.9711469933390617
$ws.Range("J6").Style = $ws.Range("G6").Style
$ws.Range("G7").Value = 0.968663215637207
$ws.Range("H7").Value = 0.970403254032135
$ws.Range("H7").Style = $ws.Range("G7").Style
$ws.Range("I7").Value = 0.9714917088534626
$ws.Range("I7").Style = $ws.Range("G7").Style
$ws.Range("J7").Value = 0.9710953235626221
$ws.Range("J7").Style = $ws.Range("G7").Style
$ws.Range("G8").Value = 0.9693394899368286
$ws.Range("H8").Value = 0.9702978730201721
$ws.Range("H8").Style = $ws.Range("G8").Style
$ws.Range("I8").Value = 0.9715176614261628
$ws.Range("I8").Style = $ws.Range("G8").Style
$ws.Range("J8").Value = 0.9710193108767271
$ws.Range("J8").Style = $ws.Range("G8").Style
$ws.Range("G9").Value = 0.9703143239021301
$ws.Range("H9").Value = 0.9701294302940369
$ws.Range("H9").Style = $ws.Range("G9").Style
$ws.Range("I9").Value = 0.971446757393018
$ws.Range("I9").Style = $ws.Range("G9").Style
$ws.Range("J9").Value = 0.9710960332304239
$ws.Range("J9").Style = $ws.Range("G9").Style
$ws.Range("G10").Value = 0.969562292098999
$ws.Range("H10").Value = 0.9702839255332947
$ws.Range("H10").Style = $ws.Range("G10").Style
$ws.Range("I10").Value = 0.9714017421985528
$ws.Range("I10").Style = $ws.Range("G10").Style
$ws.Range("J10").Value = 0.9709310252219439
$ws.Range("J10").Style = $ws.Range("G10").Style
$ws.Range("G11").Value = 0.9693045616149902
$ws.Range("H11").Value = 0.9704572558403015
$ws.Range("H11").Style = $ws.Range("G11").Style
$ws.Range("I11").Value = 0.971399820699556
$ws.Range("I11").Style = $ws.Range("G11").Style
$ws.Range("J11").Value = 0.9708862695842981
$ws.Range("J11").Style = $ws.Range("G11").Style
$ws.Range("G12").Value = 0.9697837829589844
$ws.Range("H12").Value = 0.9702953100204468
$ws.Range("H12").Style = $ws.Range("G12").Style
$ws.Range("I12").Value = 0.9713984798130514
$ws.Range("I12").Style = $ws.Range("G12").Style
$ws.Range("J12").Value = 0.9709929507225752
$ws.Range("J12").Style = $ws.Range("G12").Style
$ws.Range("G13").Value = 0.968310534954071
$ws.Range("H13").Value = 0.97027987241745
$ws.Range("H13").Style = $ws.Range("G13").Style
$ws.Range("I13").Value = 0.9714021579686458
$ws.Range("I13").Style = $ws.Range("G13").Style
$ws.Range("J13").Value = 0.9710420165210962
$ws.Range("J13").Style = $ws.Range("G13").Style
$ws.Range("G14").Value = 0.9713999032974243
$ws.Range("H14").Value = 0.9706463217735291
$ws.Range("H14").Style = $ws.Range("G14").Style
$ws.Range("I14").Value = 0.971674482704142
$ws.Range("I14").Style = $ws.Range("G14").Style
$ws.Range("J14").Value = 0.9713309016078711
$ws.Range("J14").Style = $ws.Range("G14").Style
$ws.Range("G15").Value = 0.9731541872024536
$ws.Range("H15").Value = 0.9707151651382446
$ws.Range("H15").Style = $ws.Range("G15").Style
$ws.Range("I15").Value = 0.9717204556765836
$ws.Range("I15").Style = $ws.Range("G15").Style
$ws.Range("J15").Value = 0.9713439382612705
$ws.Range("J15").Style = $ws.Range("G15").Style
$ws.Range("G16").Value = 0.9703613519668579
$ws.Range("H16").Value = 0.9704041481018066
$ws.Range("H16").Style = $ws.Range("G16").Style
$ws.Range("I16").Value = 0.9715105588283274
$ws.Range("I16").Style = $ws.Range("G16").Style
$ws.Range("J16").Value = 0.9710334651172161
$ws.Range("J16").Style = $ws.Range("G16").Style
$ws.Range("G17").Value = 0.9717307090759277
$ws.Range("H17").Value = 0.9705173373222351
$ws.Range("H17").Style = $ws.Range("G17").Style
$ws.Range("I17").Value = 0.9714847447401564
$ws.Range("I17").Style = $ws.Range("G17").Style
$ws.Range("J17").Value = 0.9709462132304907
$ws.Range("J17").Style = $ws.Range("G17").Style
$ws.Range("G18").Value = 0.9713181853294373
$ws.Range("H18").Value = 0.9703420996665955
$ws.Range("H18").Style = $ws.Range("G18").Style
$ws.Range("I18").Value = 0.9715067732832088
$ws.Range("I18").Style = $ws.Range("G18").Style
$ws.Range("J18").Value = 0.9709936566650867
$ws.Range("J18").Style = $ws.Range("G18").Style
$ws.Range("G19").Value = 0.9707207679748535
$ws.Range("H19").Value = 0.9704458117485046
$ws.Range("H19").Style = $ws.Range("G19").Style
$ws.Range("I19").Value = 0.9715038750282204
$ws.Range("I19").Style = $ws.Range("G19").Style
$ws.Range("J19").Value = 0.971008975058794
$ws.Range("J19").Style = $ws.Range("G19").Style

# --- Sheet Caso3: add Autogluon / H2O / AutoSklearn columns ---
$ws = $wb.Worksheets.Item("Caso3")

# Header row: new column titles (style matches existing unstyled G1 header cell)
$ws.Range("H1").Value = "Valori stimati autogloun"
$ws.Range("H1").Style = $ws.Range("G1").Style
$ws.Range("I1").Value = "Valori stimati h2o"
$ws.Range("I1").Style = $ws.Range("G1").Style
$ws.Range("J1").Value = "Valori stimati autosklearn"
$ws.Range("J1").Style = $ws.Range("G1").Style

$ws.Range("G2").Value = 0.9761663675308228
$ws.Range("H2").Value = 0.9725909233093262
$ws.Range("H2").Style = $ws.Range("G2").Style
$ws.Range("I2").Value = 0.9730933236636564
$ws.Range("I2").Style = $ws.Range("G2").Style
$ws.Range("J2").Value = 0.9725103303790092
$ws.Range("J2").Style = $ws.Range("G2").Style
$ws.Range("G3").Value = 0.9704825878143311
$ws.Range("H3").Value = 0.9709097146987915
$ws.Range("H3").Style = $ws.Range("G3").Style
$ws.Range("I3").Value = 0.9715747160324228
$ws.Range("I3").Style = $ws.Range("G3").Style
$ws.Range("J3").Value = 0.9711426477879286
$ws.Range("J3").Style = $ws.Range("G3").Style
$ws.Range("G4").Value = 0.9719361066818237
$ws.Range("H4").Value = 0.970124363899231
$ws.Range("H4").Style = $ws.Range("G4").Style
$ws.Range("I4").Value = 0.9712734130392966
$ws.Range("I4").Style = $ws.Range("G4").Style
$ws.Range("J4").Value = 0.9708342961966991
$ws.Range("J4").Style = $ws.Range("G4").Style
$ws.Range("G5").Value = 0.9703037142753601
$ws.Range("H5").Value = 0.9699938297271729
$ws.Range("H5").Style = $ws.Range("G5").Style
$ws.Range("I5").Value = 0.9712527984532416
$ws.Range("I5").Style = $ws.Range("G5").Style
$ws.Range("J5").Value = 0.9709512088447809
$ws.Range("J5").Style = $ws.Range("G5").Style
$ws.Range("G6").Value = 0.9703472852706909
$ws.Range("H6").Value = 0.9697715640068054
$ws.Range("H6").Style = $ws.Range("G6").Style
$ws.Range("I6").Value = 0.9712044427123936
$ws.Range("I6").Style = $ws.Range("G6").Style
$ws.Range("J6").Value = 0.9706837050616741
$ws.Range("J6").Style = $ws.Range("G6").Style
$ws.Range("G7").Value = 0.9682731032371521
$ws.Range("H7").Value = 0.970022976398468
$ws.Range("H7").Style = $ws.Range("G7").Style
$ws.Range("I7").Value = 0.9712209654304896
$ws.Range("I7").Style = $ws.Range("G7").Style
$ws.Range("J7").Value = 0.9706852622330189
$ws.Range("J7").Style = $ws.Range("G7").Style
$ws.Range("G8").Value = 0.968962550163269
$ws.Range("H8").Value = 0.9698737263679504
$ws.Range("H8").Style = $ws.Range("G8").Style
$ws.Range("I8").Value = 0.9712459327907126
$ws.Range("I8").Style = $ws.Range("G8").Style
$ws.Range("J8").Value = 0.970629608258605
$ws.Range("J8").Style = $ws.Range("G8").Style
$ws.Range("G9").Value = 0.9699459671974182
$ws.Range("H9").Value = 0.969679594039917
$ws.Range("H9").Style = $ws.Range("G9").Style
$ws.Range("I9").Value = 0.9711797904231086
$ws.Range("I9").Style = $ws.Range("G9").Style
$ws.Range("J9").Value = 0.9707021750509739
$ws.Range("J9").Style = $ws.Range("G9").Style
$ws.Range("G10").Value = 0.9691745042800903
$ws.Range("H10").Value = 0.9698253273963928
$ws.Range("H10").Style = $ws.Range("G10").Style
$ws.Range("I10").Value = 0.971136174815288
$ws.Range("I10").Style = $ws.Range("G10").Style
$ws.Range("J10").Value = 0.9705584030598402
$ws.Range("J10").Style = $ws.Range("G10").Style
$ws.Range("G11").Value = 0.968936562538147
$ws.Range("H11").Value = 0.9701217412948608
$ws.Range("H11").Style = $ws.Range("G11").Style
$ws.Range("I11").Value = 0.9711325294274888
$ws.Range("I11").Style = $ws.Range("G11").Style
$ws.Range("J11").Value = 0.9704122822731733
$ws.Range("J11").Style = $ws.Range("G11").Style
$ws.Range("G12").Value = 0.969379723072052
$ws.Range("H12").Value = 0.9699321985244751
$ws.Range("H12").Style = $ws.Range("G12").Style
$ws.Range("I12").Value = 0.9711326756830267
$ws.Range("I12").Style = $ws.Range("G12").Style
$ws.Range("J12").Value = 0.9705849774181843
$ws.Range("J12").Style = $ws.Range("G12").Style
$ws.Range("G13").Value = 0.9679166674613953
$ws.Range("H13").Value = 0.9698776006698608
$ws.Range("H13").Style = $ws.Range("G13").Style
$ws.Range("I13").Value = 0.9711338085195156
$ws.Range("I13").Style = $ws.Range("G13").Style
$ws.Range("J13").Value = 0.9706240314990282
$ws.Range("J13").Style = $ws.Range("G13").Style
$ws.Range("G14").Value = 0.9710345268249512
$ws.Range("H14").Value = 0.9702818989753723
$ws.Range("H14").Style = $ws.Range("G14").Style
$ws.Range("I14").Value = 0.9714042566292288
$ws.Range("I14").Style = $ws.Range("G14").Style
$ws.Range("J14").Value = 0.9708847478032112
$ws.Range("J14").Style = $ws.Range("G14").Style
$ws.Range("G15").Value = 0.972794771194458
$ws.Range("H15").Value = 0.970242440700531
$ws.Range("H15").Style = $ws.Range("G15").Style
$ws.Range("I15").Value = 0.9714490114506836
$ws.Range("I15").Style = $ws.Range("G15").Style
$ws.Range("J15").Value = 0.9708912093192339
$ws.Range("J15").Style = $ws.Range("G15").Style
$ws.Range("G16").Value = 0.9699744582176208
$ws.Range("H16").Value = 0.9700605273246765
$ws.Range("H16").Style = $ws.Range("G16").Style
$ws.Range("I16").Value = 0.971242152179398
$ws.Range("I16").Style = $ws.Range("G16").Style
$ws.Range("J16").Value = 0.9706978052854538
$ws.Range("J16").Style = $ws.Range("G16").Style
$ws.Range("G17").Value = 0.9713557958602905
$ws.Range("H17").Value = 0.9701403379440308
$ws.Range("H17").Style = $ws.Range("G17").Style
$ws.Range("I17").Value = 0.971218490865024
$ws.Range("I17").Style = $ws.Range("G17").Style
$ws.Range("J17").Value = 0.9706005435436964
$ws.Range("J17").Style = $ws.Range("G17").Style
$ws.Range("G18").Value = 0.9709411859512329
$ws.Range("H18").Value = 0.9700134992599487
$ws.Range("H18").Style = $ws.Range("G18").Style
$ws.Range("I18").Value = 0.9712422761025408
$ws.Range("I18").Style = $ws.Range("G18").Style
$ws.Range("J18").Value = 0.9706536922603846
$ws.Range("J18").Style = $ws.Range("G18").Style
$ws.Range("G19").Value = 0.9703565835952759
$ws.Range("H19").Value = 0.9699594974517822
$ws.Range("H19").Style = $ws.Range("G19").Style
$ws.Range("I19").Value = 0.9712377303692026
$ws.Range("I19").Style = $ws.Range("G19").Style
$ws.Range("J19").Value = 0.9707414116710424
$ws.Range("J19").Style = $ws.Range("G19").Style

